$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update StartDate/EndDate for the first four reservations (HR001-HR004, rows 2-5)
# from 10/02/2019 - 12/02/2019 to 04/02/2019 - 06/02/2019
$ws.Cells.Item(2, 3).Value = "04/02/2019"
$ws.Cells.Item(2, 4).Value = "06/02/2019"
$ws.Cells.Item(3, 3).Value = "04/02/2019"
$ws.Cells.Item(3, 4).Value = "06/02/2019"
$ws.Cells.Item(4, 3).Value = "04/02/2019"
$ws.Cells.Item(4, 4).Value = "06/02/2019"
$ws.Cells.Item(5, 3).Value = "04/02/2019"
$ws.Cells.Item(5, 4).Value = "06/02/2019"

# Update the sheet view: scroll/select column B instead of A/E (topLeftCell B1, active cell B10)
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B10").Select()
